# Updates the LR-pairs NATMI TPM metrics (columns E:T, rows 2-10) with recomputed values,
# matching the "update scripts wuth new tpm" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.5586043333333334
$ws.Cells.Item(2,8).Value = 1.675813
$ws.Cells.Item(2,9).Value = 0.01643366487114074
$ws.Cells.Item(2,10).Value = 0.01643366487114074
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.09554499999999999
$ws.Cells.Item(2,14).Value = 0.286635
$ws.Cells.Item(2,15).Value = 0.01029975823317688
$ws.Cells.Item(2,16).Value = 0.01029975823317688
$ws.Cells.Item(2,17).Value = 0.05337185102833333
$ws.Cells.Item(2,18).Value = 0.480346659255
$ws.Cells.Item(2,19).Value = 0.0001692627750578015
$ws.Cells.Item(2,20).Value = 0.0001692627750578015

# Row 3
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.5586043333333334
$ws.Cells.Item(3,8).Value = 1.675813
$ws.Cells.Item(3,9).Value = 0.01643366487114074
$ws.Cells.Item(3,10).Value = 0.01643366487114074
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.721182333333333
$ws.Cells.Item(3,14).Value = 11.163547
$ws.Cells.Item(3,15).Value = 0.4011437372432085
$ws.Cells.Item(3,16).Value = 0.4011437372432086
$ws.Cells.Item(3,17).Value = 2.078668576523445
$ws.Cells.Item(3,18).Value = 18.708017188711
$ws.Cells.Item(3,19).Value = 0.006592261743011827
$ws.Cells.Item(3,20).Value = 0.006592261743011827

# Row 4
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.5586043333333334
$ws.Cells.Item(4,8).Value = 1.675813
$ws.Cells.Item(4,9).Value = 0.01643366487114074
$ws.Cells.Item(4,10).Value = 0.01643366487114074
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.459703999999999
$ws.Cells.Item(4,14).Value = 16.379112
$ws.Cells.Item(4,15).Value = 0.5885565045236145
$ws.Cells.Item(4,16).Value = 0.5885565045236146
$ws.Cells.Item(4,17).Value = 3.049814313117333
$ws.Cells.Item(4,18).Value = 27.448328818056
$ws.Cells.Item(4,19).Value = 0.00967214035307111
$ws.Cells.Item(4,20).Value = 0.00967214035307111

# Row 5
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 21.275312
$ws.Cells.Item(5,8).Value = 63.825936
$ws.Cells.Item(5,9).Value = 0.6259016025719319
$ws.Cells.Item(5,10).Value = 0.6259016025719319
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.09554499999999999
$ws.Cells.Item(5,14).Value = 0.286635
$ws.Cells.Item(5,15).Value = 0.01029975823317688
$ws.Cells.Item(5,16).Value = 0.01029975823317688
$ws.Cells.Item(5,17).Value = 2.03274968504
$ws.Cells.Item(5,18).Value = 18.29474716536
$ws.Cells.Item(5,19).Value = 0.00644663518424886
$ws.Cells.Item(5,20).Value = 0.00644663518424886

# Row 6
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 21.275312
$ws.Cells.Item(6,8).Value = 63.825936
$ws.Cells.Item(6,9).Value = 0.6259016025719319
$ws.Cells.Item(6,10).Value = 0.6259016025719319
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.721182333333333
$ws.Cells.Item(6,14).Value = 11.163547
$ws.Cells.Item(6,15).Value = 0.4011437372432085
$ws.Cells.Item(6,16).Value = 0.4011437372432086
$ws.Cells.Item(6,17).Value = 79.16931515055467
$ws.Cells.Item(6,18).Value = 712.523836354992
$ws.Cells.Item(6,19).Value = 0.2510765080022182
$ws.Cells.Item(6,20).Value = 0.2510765080022182

# Row 7
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 21.275312
$ws.Cells.Item(7,8).Value = 63.825936
$ws.Cells.Item(7,9).Value = 0.6259016025719319
$ws.Cells.Item(7,10).Value = 0.6259016025719319
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.459703999999999
$ws.Cells.Item(7,14).Value = 16.379112
$ws.Cells.Item(7,15).Value = 0.5885565045236145
$ws.Cells.Item(7,16).Value = 0.5885565045236146
$ws.Cells.Item(7,17).Value = 116.156906027648
$ws.Cells.Item(7,18).Value = 1045.412154248832
$ws.Cells.Item(7,19).Value = 0.3683784593854648
$ws.Cells.Item(7,20).Value = 0.3683784593854649

# Row 8
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 12.157548
$ws.Cells.Item(8,8).Value = 36.472644
$ws.Cells.Item(8,9).Value = 0.3576647325569273
$ws.Cells.Item(8,10).Value = 0.3576647325569273
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.09554499999999999
$ws.Cells.Item(8,14).Value = 0.286635
$ws.Cells.Item(8,15).Value = 0.01029975823317688
$ws.Cells.Item(8,16).Value = 0.01029975823317688
$ws.Cells.Item(8,17).Value = 1.16159292366
$ws.Cells.Item(8,18).Value = 10.45433631294
$ws.Cells.Item(8,19).Value = 0.00368386027387022
$ws.Cells.Item(8,20).Value = 0.00368386027387022

# Row 9
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 12.157548
$ws.Cells.Item(9,8).Value = 36.472644
$ws.Cells.Item(9,9).Value = 0.3576647325569273
$ws.Cells.Item(9,10).Value = 0.3576647325569273
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.721182333333333
$ws.Cells.Item(9,14).Value = 11.163547
$ws.Cells.Item(9,15).Value = 0.4011437372432085
$ws.Cells.Item(9,16).Value = 0.4011437372432086
$ws.Cells.Item(9,17).Value = 45.240452834252
$ws.Cells.Item(9,18).Value = 407.164075508268
$ws.Cells.Item(9,19).Value = 0.1434749674979785
$ws.Cells.Item(9,20).Value = 0.1434749674979785

# Row 10
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 12.157548
$ws.Cells.Item(10,8).Value = 36.472644
$ws.Cells.Item(10,9).Value = 0.3576647325569273
$ws.Cells.Item(10,10).Value = 0.3576647325569273
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.459703999999999
$ws.Cells.Item(10,14).Value = 16.379112
$ws.Cells.Item(10,15).Value = 0.5885565045236145
$ws.Cells.Item(10,16).Value = 0.5885565045236146
$ws.Cells.Item(10,17).Value = 66.37661344579199
$ws.Cells.Item(10,18).Value = 597.3895210121281
$ws.Cells.Item(10,19).Value = 0.2105059047850786
$ws.Cells.Item(10,20).Value = 0.2105059047850786

